$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.872588634490967
$ws.Range("B1").Value = 5.900938987731934
$ws.Range("C1").Value = 6.616466522216797
$ws.Range("D1").Value = 9.912729263305664
$ws.Range("E1").Value = 5.800085067749023
